# "print stacked bar chart with error bars"
# sheet_two gets the same Temperature2/Temperature3/Average/Deviation
# columns (C:F) that sheet_one already has, with a fresh set of sample
# values in column B, and becomes the active tab/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("sheet_one")
$ws2 = $wb.Worksheets.Item("sheet_two")

# --- sheet_two: new data -------------------------------------------------

# Headers for the newly added columns (reuses existing shared strings).
$ws2.Range("C1").Value = "Temperature2"
$ws2.Range("D1").Value = "Temperature3"
$ws2.Range("E1").Value = "Average"
$ws2.Range("F1").Value = "Deviation"

# Column B: new sample readings (descending), replacing the old series.
$bValues = @(30,29,28,27,26,25,24,23,22,21,20,19,18,17,16,15,14,13,12,11,10,9,8,7,6)
for ($i = 0; $i -lt $bValues.Count; $i++) {
    $ws2.Cells.Item($i + 2, 2).Value = $bValues[$i]
}

# Column C: second temperature series.
$cValues = @(20,21,20,21,20,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4,20.4)
for ($i = 0; $i -lt $cValues.Count; $i++) {
    $ws2.Cells.Item($i + 2, 3).Value = $cValues[$i]
}

# Column D: Temperature3 = B * 0.9 (row 2 alone, then a shared formula for 3:26).
$ws2.Range("D2").Formula = "=B2*0.9"
$ws2.Range("D3:D26").Formula = "=B3*0.9"

# Column E: Average of B:D.
$ws2.Range("E2").Formula = "=AVERAGE(B2:D2)"
$ws2.Range("E3:E26").Formula = "=AVERAGE(B3:D3)"

# Column F: sample deviation of B:D.
$ws2.Range("F2").Formula = "=STDEVA(B2:D2)"
$ws2.Range("F3:F26").Formula = "=STDEVA(B3:D3)"

# New columns C:D sized like sheet_one's equivalent auto-fit columns.
$ws2.Range("C1:D1").ColumnWidth = 13.140625

# --- selections / active sheet -------------------------------------------

# sheet_one loses its selection/tab focus, keeping only a moved selection.
$ws1.Activate()
$ws1.Range("E31").Select()

# sheet_two becomes the tab that is selected/active in the saved workbook.
$ws2.Activate()
$ws2.Range("F30").Select()
